$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "rxxx"
$ws.Range("B9").Value = "greg"
$ws.Range("C9").Value = "temp profile off on setup"
$ws.Range("D9").Value = "2025-09-30 13:25:50"
